# dlgMotorcycleInsurance.xlsx bugfix: add per-section XPath locators for the
# insurance-form web page and select the newly written header range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "//*[@id='insurance-form']/div/section[1]"
$ws.Range("D1").Value = "//*[@id='insurance-form']/div/section[2]"
$ws.Range("E1").Value = "//*[@id='insurance-form']/div/section[3]"
$ws.Range("F1").Value = "//*[@id='insurance-form']/div/section[4]"
$ws.Range("G1").Value = "//*[@id='insurance-form']/div/section[5]"

$ws.Range("C1:G1").Select()
